$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.105.34"
$ws.Range("E2").Value = "  +5.14%  "
$ws.Range("D3").Value = "3.370.19"
$ws.Range("E3").Value = "  +5.73%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'572.47"
$ws.Range("E5").Value = "  +7.03%  "
$ws.Range("D6").Value = "'152.66"
$ws.Range("E6").Value = "  +5.05%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").Value = "3.380.42"
$ws.Range("E8").Value = "  +5.69%  "
$ws.Range("D9").Value = "'0.527"
$ws.Range("E9").Value = "  +0.16%  "
$ws.Range("E10").Value = "  +1.60%  "
$ws.Range("E11").Value = "  +5.92%  "
$ws.Range("E12").Value = "  +1.55%  "
$ws.Range("D13").Value = "3.955.42"
$ws.Range("E13").Value = "  +5.92%  "
$ws.Range("E14").Value = "  +0.25%  "
$ws.Range("D15").Value = "'26.98"
$ws.Range("E15").Value = "  +4.37%  "
$ws.Range("D16").Value = "'0.0000180"
$ws.Range("E16").Value = "  +4.48%  "
$ws.Range("D17").Value = "63.109.06"
$ws.Range("E17").Value = "  +5.11%  "
$ws.Range("D18").Value = "3.356.41"
$ws.Range("E18").Value = "  +5.25%  "
$ws.Range("E19").Value = "  +1.22%  "
$ws.Range("D20").Value = "'13.91"
$ws.Range("E20").Value = "  +5.07%  "
$ws.Range("D21").Value = "'8.41"
$ws.Range("E21").Value = "  +2.79%  "
$ws.Range("D22").Value = "'385.80"
$ws.Range("E22").Value = "  +4.58%  "
$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = "  -0.19%  "
$ws.Range("D24").Value = "'0.533"
$ws.Range("E24").Value = "  +2.47%  "
$ws.Range("D25").Value = "'70.38"
$ws.Range("E25").Value = "  +1.21%  "
$ws.Range("E26").Value = "  +6.52%  "
$ws.Range("D27").Value = "'9.23"
$ws.Range("E27").Value = "  +6.64%  "
$ws.Range("D28").Value = "0.0₃0970"
$ws.Range("E28").Value = "  +11.99%  "
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("E31").Value = "  +3.45%  "
$ws.Range("D32").Value = "'1.31"
$ws.Range("E32").Value = "  +10.61%  "
$ws.Range("E33").Value = "  +5.36%  "
$ws.Range("D34").Value = "'6.29"
$ws.Range("E34").Value = "  +4.12%  "
$ws.Range("E35").Value = "  +2.21%  "
$ws.Range("E36").Value = "  +9.31%  "
$ws.Range("D37").Value = "'157.97"
$ws.Range("E37").Value = "  +1.21%  "
$ws.Range("E38").Value = "  +12.31%  "
$ws.Range("D39").Value = "'27.38"
$ws.Range("E39").Value = "  +5.10%  "
$ws.Range("D40").Value = "2.889.31"
$ws.Range("E40").Value = "  +2.64%  "
$ws.Range("E41").Value = "  +10.51%  "
$ws.Range("D42").Value = "'0.0742"
$ws.Range("E42").Value = "  +5.53%  "
$ws.Range("D43").Value = "'40.82"
$ws.Range("E43").Value = "  +2.97%  "
$ws.Range("D44").Value = "'0.748"
$ws.Range("E44").Value = "  +4.44%  "
$ws.Range("D45").Value = "'4.24"
$ws.Range("E45").Value = "  +0.57%  "
$ws.Range("E46").Value = "  +5.77%  "
$ws.Range("D47").Value = "3.421.18"
$ws.Range("E47").Value = "  +5.96%  "
$ws.Range("E48").Value = "  +6.44%  "
$ws.Range("D49").Value = "'300.73"
$ws.Range("E49").Value = "  +13.75%  "
$ws.Range("E50").Value = "  -1.64%  "
$ws.Range("E51").Value = "  +2.46%  "
